# dt_full_qoq_PRIVCON_AVERAGE_1_9 : extend the naive-forecast table by one more
# period (new column BB / new row 83), per "Included EQUIPMENT eval, updated
# folder structure, allowed for multiple archive excels".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header: append the next quarter-end date in BB1, matching BA1 formatting ---
$ws.Range("BA1").Copy()
$ws.Range("BB1").PasteSpecial(-4122)
$ws.Range("BB1").Value = 45986

# --- Rows 2-70: column BB repeats the last known (BA) forecast value for that row ---
$ws.Range("BB2").Value = -0.2522510312815029
$ws.Range("BB3").Value = 0.8253228627186644
$ws.Range("BB4").Value = -0.8956080111821194
$ws.Range("BB5").Value = 0.6802092152583725
$ws.Range("BB6").Value = -0.2702467235064887
$ws.Range("BB7").Value = 0.2360788655137185
$ws.Range("BB8").Value = 0.4710547017577227
$ws.Range("BB9").Value = -0.6739669439274962
$ws.Range("BB10").Value = -0.2065125907403171
$ws.Range("BB11").Value = 0.1861693474447037
$ws.Range("BB12").Value = 0.7530542497296722
$ws.Range("BB13").Value = 0.3106180987640244
$ws.Range("BB14").Value = 0.5709322285556908
$ws.Range("BB15").Value = 1.247003875094862
$ws.Range("BB16").Value = -0.776883278166693
$ws.Range("BB17").Value = 1.384511819460911
$ws.Range("BB18").Value = 0.09417825394891111
$ws.Range("BB19").Value = -0.1288380179123294
$ws.Range("BB20").Value = 0.1460429756367461
$ws.Range("BB21").Value = -0.0194429241770564
$ws.Range("BB22").Value = 0.2041957853813301
$ws.Range("BB23").Value = 0.1772414135220401
$ws.Range("BB24").Value = 0.7756459370471021
$ws.Range("BB25").Value = 0.5676350012738425
$ws.Range("BB26").Value = -0.7653316860800885
$ws.Range("BB27").Value = 0.322227538137227
$ws.Range("BB28").Value = -0.0579621925135001
$ws.Range("BB29").Value = 0.6475935709367775
$ws.Range("BB30").Value = 0.9603367340567104
$ws.Range("BB31").Value = 0.3700548251239582
$ws.Range("BB32").Value = 0.4443178743943008
$ws.Range("BB33").Value = 0.5552956277764309
$ws.Range("BB34").Value = 0.3556667645181193
$ws.Range("BB35").Value = 0.6924895145077272
$ws.Range("BB36").Value = 0.3519868436780342
$ws.Range("BB37").Value = 0.489196258618918
$ws.Range("BB38").Value = 0.1561519231779869
$ws.Range("BB39").Value = 0.5297299217112936
$ws.Range("BB40").Value = 0.7664191671286744
$ws.Range("BB41").Value = 0.04527341468607915
$ws.Range("BB42").Value = 0.06335359735614077
$ws.Range("BB43").Value = 0.08598758370690973
$ws.Range("BB44").Value = 0.2196343350075409
$ws.Range("BB45").Value = 0.03810891122928695
$ws.Range("BB46").Value = 0.4286160255108911
$ws.Range("BB47").Value = 0.9
$ws.Range("BB48").Value = 0.2
$ws.Range("BB49").Value = 0.2
$ws.Range("BB50").Value = 0
$ws.Range("BB51").Value = -2.213339122522456
$ws.Range("BB52").Value = -11.49785608241407
$ws.Range("BB53").Value = 11.22930999924247
$ws.Range("BB54").Value = -2.250986781122748
$ws.Range("BB55").Value = -5.40098554941693
$ws.Range("BB56").Value = 3.942549781810342
$ws.Range("BB57").Value = 6.284666192508709
$ws.Range("BB58").Value = -1.347757551663406
$ws.Range("BB59").Value = 1.396500038188336
$ws.Range("BB60").Value = 0.1885473380929312
$ws.Range("BB61").Value = 1.224073604180177
$ws.Range("BB62").Value = -1.674179157827311
$ws.Range("BB63").Value = -0.6601113848982436
$ws.Range("BB64").Value = 0.2658000717656392
$ws.Range("BB65").Value = -0.01892863903084674
$ws.Range("BB66").Value = 0.4261262404008619
$ws.Range("BB67").Value = -0.07433428650158191
$ws.Range("BB68").Value = 0.09298591595782568
$ws.Range("BB69").Value = 0.2693878213604393
$ws.Range("BB70").Value = 0.2408767182737677

# --- Rows 71-82: column BB receives the newly computed forecast edge values ---
$ws.Range("BB71").Value = 0.5597354586130052
$ws.Range("BB72").Value = 0.1368731201391853
$ws.Range("BB73").Value = -0.2551464291630765
$ws.Range("BB74").Value = -0.2551464291630765
$ws.Range("BB75").Value = -0.2551464291630765
$ws.Range("BB76").Value = -0.2551464291630765
$ws.Range("BB77").Value = -0.2551464291630765
$ws.Range("BB78").Value = -0.2551464291630765
$ws.Range("BB79").Value = -0.2551464291630765
$ws.Range("BB80").Value = -0.2551464291630765
$ws.Range("BB81").Value = -0.2551464291630765
$ws.Range("BB82").Value = -0.2551464291630765

# --- New row 83: next period date in column A (formatted like A82), plus its BB forecast ---
$ws.Range("A82").Copy()
$ws.Range("A83").PasteSpecial(-4122)
$ws.Range("A83").Value = 46934
$ws.Range("BB83").Value = -0.2551464291630765

$excel.CutCopyMode = $false

Write-Host "done"
